$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions data pull): update Price (D) and
# Volume(1h) (E) columns. Some Price values are plain decimals that Excel
# would otherwise auto-convert to numbers, so those cells are forced to
# Text format first to preserve the original text storage.

$ws.Range("D2").Value = "61.723.43"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").Value = "2.573.51"
$ws.Range("E3").Value = "  -4.75%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.86"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.11"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "3.027.74"
$ws.Range("E13").Value = "  -5.02%  "
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").Value = "61.604.85"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").Value = "2.576.92"
$ws.Range("E17").Value = "  -5.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.50"
$ws.Range("E18").Value = "  -4.83%  "
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.43"
$ws.Range("E20").Value = "  -3.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  -4.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.491"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.18"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.05"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.48"
$ws.Range("E28").Value = "  +5.66%  "
$ws.Range("D29").Value = "0.0₃0837"
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.55"
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.72"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.13"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "331.75"
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.917"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.87"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.93"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.50"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.59"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "2.122.88"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.54"
$ws.Range("E48").Value = "  -4.30%  "
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("E51").Value = "  -1.73%  "
